$wb = $excel.ActiveWorkbook

# --- top_products: update category names + quantities (rows 2-11) ---
$ws = $wb.Worksheets.Item("top_products")
$ws.Cells.Item(2, 1).Value = "Mutfak Eşyası"
$ws.Cells.Item(2, 2).Value = 37
$ws.Cells.Item(3, 1).Value = "Sulama Sistemi"
$ws.Cells.Item(3, 2).Value = 35
$ws.Cells.Item(4, 1).Value = "Araç Aksesuarı"
$ws.Cells.Item(4, 2).Value = 34
$ws.Cells.Item(5, 1).Value = "Bahçe Aletleri"
$ws.Cells.Item(5, 2).Value = 34
$ws.Cells.Item(6, 1).Value = "Dekorasyon"
$ws.Cells.Item(6, 2).Value = 34
$ws.Cells.Item(7, 1).Value = "Klavye"
$ws.Cells.Item(7, 2).Value = 33
$ws.Cells.Item(8, 1).Value = "Bisiklet"
$ws.Cells.Item(8, 2).Value = 33
$ws.Cells.Item(9, 1).Value = "Çanta"
$ws.Cells.Item(9, 2).Value = 33
$ws.Cells.Item(10, 1).Value = "Tenis Raketi"
$ws.Cells.Item(10, 2).Value = 33
$ws.Cells.Item(11, 1).Value = "Akademik"
$ws.Cells.Item(11, 2).Value = 32

# --- category_prices: update mean/min/max (rows 2-11); category labels unchanged ---
$ws = $wb.Worksheets.Item("category_prices")
$ws.Cells.Item(2, 2).Value = 893.5477192982456
$ws.Cells.Item(2, 3).Value = 103.34
$ws.Cells.Item(2, 4).Value = 1995.58
$ws.Cells.Item(3, 2).Value = 2774.251754385965
$ws.Cells.Item(3, 3).Value = 526.98
$ws.Cells.Item(3, 4).Value = 4993.84
$ws.Cells.Item(4, 2).Value = 1435.741178369653
$ws.Cells.Item(4, 3).Value = 132.78
$ws.Cells.Item(4, 4).Value = 2999.64
$ws.Cells.Item(5, 2).Value = 237.8455357142857
$ws.Cells.Item(5, 3).Value = 56.6
$ws.Cells.Item(5, 4).Value = 473.58
$ws.Cells.Item(6, 2).Value = 87.4574358974359
$ws.Cells.Item(6, 3).Value = 22.7
$ws.Cells.Item(6, 4).Value = 154.85
$ws.Cells.Item(7, 2).Value = 178.5532075471698
$ws.Cells.Item(7, 3).Value = 32.93
$ws.Cells.Item(7, 4).Value = 299.45
$ws.Cells.Item(8, 2).Value = 46.61862068965517
$ws.Cells.Item(8, 3).Value = 6.21
$ws.Cells.Item(8, 4).Value = 99.77
$ws.Cells.Item(9, 2).Value = 562.9056923076923
$ws.Cells.Item(9, 3).Value = 77.73
$ws.Cells.Item(9, 4).Value = 998.67
$ws.Cells.Item(10, 2).Value = 260.8737391304348
$ws.Cells.Item(10, 3).Value = 63.71
$ws.Cells.Item(10, 4).Value = 514.8
$ws.Cells.Item(11, 2).Value = 1005.694345238095
$ws.Cells.Item(11, 3).Value = 171.02
$ws.Cells.Item(11, 4).Value = 1974.15

# --- recent_sales: replace data rows with new (longer) list (rows 2-40) ---
$ws = $wb.Worksheets.Item("recent_sales")
$fmtSrc = $ws.Range("A2:B2")
$fmtDst = $ws.Range("A2:B40")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(2, 1).Value = "Biyografi"
$ws.Cells.Item(2, 2).Value = 11
$ws.Cells.Item(3, 1).Value = "Koşu Ayakkabısı"
$ws.Cells.Item(3, 2).Value = 8
$ws.Cells.Item(4, 1).Value = "Mouse"
$ws.Cells.Item(4, 2).Value = 8
$ws.Cells.Item(5, 1).Value = "Tenis Raketi"
$ws.Cells.Item(5, 2).Value = 8
$ws.Cells.Item(6, 1).Value = "Ev Tekstili"
$ws.Cells.Item(6, 2).Value = 7
$ws.Cells.Item(7, 1).Value = "Çanta"
$ws.Cells.Item(7, 2).Value = 7
$ws.Cells.Item(8, 1).Value = "Ajanda"
$ws.Cells.Item(8, 2).Value = 6
$ws.Cells.Item(9, 1).Value = "Defter"
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(10, 1).Value = "Tablet"
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(11, 1).Value = "Yoga Matı"
$ws.Cells.Item(11, 2).Value = 5
$ws.Cells.Item(12, 1).Value = "Puzzle"
$ws.Cells.Item(12, 2).Value = 5
$ws.Cells.Item(13, 1).Value = "Bisiklet"
$ws.Cells.Item(13, 2).Value = 5
$ws.Cells.Item(14, 1).Value = "Sulama Sistemi"
$ws.Cells.Item(14, 2).Value = 5
$ws.Cells.Item(15, 1).Value = "Gömlek"
$ws.Cells.Item(15, 2).Value = 5
$ws.Cells.Item(16, 1).Value = "Kişisel Gelişim"
$ws.Cells.Item(16, 2).Value = 5
$ws.Cells.Item(17, 1).Value = "Perde"
$ws.Cells.Item(17, 2).Value = 4
$ws.Cells.Item(18, 1).Value = "Pantolon"
$ws.Cells.Item(18, 2).Value = 4
$ws.Cells.Item(19, 1).Value = "Futbol Topu"
$ws.Cells.Item(19, 2).Value = 4
$ws.Cells.Item(20, 1).Value = "Araç Aksesuarı"
$ws.Cells.Item(20, 2).Value = 4
$ws.Cells.Item(21, 1).Value = "Monitör"
$ws.Cells.Item(21, 2).Value = 4
$ws.Cells.Item(22, 1).Value = "Klavye"
$ws.Cells.Item(22, 2).Value = 4
$ws.Cells.Item(23, 1).Value = "Elbise"
$ws.Cells.Item(23, 2).Value = 4
$ws.Cells.Item(24, 1).Value = "Akıllı Telefon"
$ws.Cells.Item(24, 2).Value = 3
$ws.Cells.Item(25, 1).Value = "Şampuan"
$ws.Cells.Item(25, 2).Value = 3
$ws.Cells.Item(26, 1).Value = "Bahçe Aletleri"
$ws.Cells.Item(26, 2).Value = 3
$ws.Cells.Item(27, 1).Value = "Oyun Konsolu"
$ws.Cells.Item(27, 2).Value = 3
$ws.Cells.Item(28, 1).Value = "Oto Parfümü"
$ws.Cells.Item(28, 2).Value = 3
$ws.Cells.Item(29, 1).Value = "Kulaklık"
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(30, 1).Value = "Güneş Kremi"
$ws.Cells.Item(30, 2).Value = 3
$ws.Cells.Item(31, 1).Value = "Parfüm"
$ws.Cells.Item(31, 2).Value = 3
$ws.Cells.Item(32, 1).Value = "Laptop"
$ws.Cells.Item(32, 2).Value = 3
$ws.Cells.Item(33, 1).Value = "Şapka"
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(34, 1).Value = "Tarih"
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(35, 1).Value = "Makyaj Seti"
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(36, 1).Value = "Uzaktan Kumandalı Araba"
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(37, 1).Value = "Bitki"
$ws.Cells.Item(37, 2).Value = 1
$ws.Cells.Item(38, 1).Value = "Boya Seti"
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(39, 1).Value = "Dergi"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(40, 1).Value = "Yağ"
$ws.Cells.Item(40, 2).Value = 1

# --- payment_analysis: update transaction counts + totals (rows 2-6) ---
$ws = $wb.Worksheets.Item("payment_analysis")
$ws.Cells.Item(2, 2).Value = 128
$ws.Cells.Item(2, 3).Value = 279417.5514285714
$ws.Cells.Item(3, 2).Value = 119
$ws.Cells.Item(3, 3).Value = 321983.48
$ws.Cells.Item(4, 2).Value = 118
$ws.Cells.Item(4, 3).Value = 225420.4833333333
$ws.Cells.Item(5, 2).Value = 111
$ws.Cells.Item(5, 3).Value = 336205.8333333334
$ws.Cells.Item(6, 2).Value = 119
$ws.Cells.Item(6, 3).Value = 247446.3386666666

# --- missing_data_report: reorder fields + update counts/pcts (rows 2-12) ---
$ws = $wb.Worksheets.Item("missing_data_report")
$ws.Cells.Item(2, 1).Value = "satisfaction_score"
$ws.Cells.Item(2, 2).Value = 103
$ws.Cells.Item(2, 3).Value = 17.3109243697479
$ws.Cells.Item(3, 1).Value = "purchase_date"
$ws.Cells.Item(3, 2).Value = 8
$ws.Cells.Item(3, 3).Value = 1.344537815126051
$ws.Cells.Item(4, 1).Value = "price"
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 1.008403361344538
$ws.Cells.Item(5, 1).Value = "product_name"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(6, 1).Value = "customer_id"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(7, 1).Value = "id"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(8, 1).Value = "category"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 1).Value = "quantity"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(10, 1).Value = "payment_method"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 1).Value = "shipping_cost"
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(12, 1).Value = "discount_applied"
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0

# --- cleaning_comparison: update before/after/filled counts (rows 2-12) ---
$ws = $wb.Worksheets.Item("cleaning_comparison")
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(6, 2).Value = 6
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 6
$ws.Cells.Item(7, 2).Value = 8
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 8
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(9, 2).Value = 103
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 103
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0

# --- price_updates: replace rows with full recalculated price list (rows 2-44) ---
$ws = $wb.Worksheets.Item("price_updates")
$fmtSrc = $ws.Range("A2:C2")
$fmtDst = $ws.Range("A2:C44")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Akıllı Telefon"
$ws.Cells.Item(2, 3).Value = 3051.68
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Kamera"
$ws.Cells.Item(3, 3).Value = 2775.02
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Klavye"
$ws.Cells.Item(4, 3).Value = 3051.68
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Kulaklık"
$ws.Cells.Item(5, 3).Value = 2496.83
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Laptop"
$ws.Cells.Item(6, 3).Value = 2496.83
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Tablet"
$ws.Cells.Item(7, 3).Value = 3051.68
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "Ayakkabı"
$ws.Cells.Item(8, 3).Value = 261.63
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Ceket"
$ws.Cells.Item(9, 3).Value = 214.06
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gömlek"
$ws.Cells.Item(10, 3).Value = 214.06
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "T-shirt"
$ws.Cells.Item(11, 3).Value = 214.06
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Çanta"
$ws.Cells.Item(12, 3).Value = 214.06
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Şapka"
$ws.Cells.Item(13, 3).Value = 261.63
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Bilim Kurgu"
$ws.Cells.Item(14, 3).Value = 96.2
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Biyografi"
$ws.Cells.Item(15, 3).Value = 96.2
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Kişisel Gelişim"
$ws.Cells.Item(16, 3).Value = 78.71
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "Roman"
$ws.Cells.Item(17, 3).Value = 78.71
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "Tarih"
$ws.Cells.Item(18, 3).Value = 78.71
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "Çocuk Kitabı"
$ws.Cells.Item(19, 3).Value = 96.2
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Bisiklet"
$ws.Cells.Item(20, 3).Value = 905.12
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Spor Çantası"
$ws.Cells.Item(21, 3).Value = 1106.26
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "Tenis Raketi"
$ws.Cells.Item(22, 3).Value = 905.12
$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "Yüzme Gözlüğü"
$ws.Cells.Item(23, 3).Value = 1106.26
$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "Aydınlatma"
$ws.Cells.Item(24, 3).Value = 1579.32
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "Dekorasyon"
$ws.Cells.Item(25, 3).Value = 1292.17
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "Halı"
$ws.Cells.Item(26, 3).Value = 1292.17
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "Mutfak Eşyası"
$ws.Cells.Item(27, 3).Value = 1579.32
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "Perde"
$ws.Cells.Item(28, 3).Value = 1292.17
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "Krem"
$ws.Cells.Item(29, 3).Value = 160.7
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "Parfüm"
$ws.Cells.Item(30, 3).Value = 196.41
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "Ruj"
$ws.Cells.Item(31, 3).Value = 160.7
$ws.Cells.Item(32, 1).Value = 30
$ws.Cells.Item(32, 2).Value = "Şampuan"
$ws.Cells.Item(32, 3).Value = 196.41
$ws.Cells.Item(33, 1).Value = 31
$ws.Cells.Item(33, 2).Value = "Bebek"
$ws.Cells.Item(33, 3).Value = 234.79
$ws.Cells.Item(34, 1).Value = 32
$ws.Cells.Item(34, 2).Value = "Eğitici Oyuncak"
$ws.Cells.Item(34, 3).Value = 286.96
$ws.Cells.Item(35, 1).Value = 33
$ws.Cells.Item(35, 2).Value = "Lego"
$ws.Cells.Item(35, 3).Value = 234.79
$ws.Cells.Item(36, 1).Value = 34
$ws.Cells.Item(36, 2).Value = "Peluş Oyuncak"
$ws.Cells.Item(36, 3).Value = 234.79
$ws.Cells.Item(37, 1).Value = 35
$ws.Cells.Item(37, 2).Value = "Bahçe Aletleri"
$ws.Cells.Item(37, 3).Value = 804.19
$ws.Cells.Item(38, 1).Value = 36
$ws.Cells.Item(38, 2).Value = "Gübre"
$ws.Cells.Item(38, 3).Value = 982.9
$ws.Cells.Item(39, 1).Value = 37
$ws.Cells.Item(39, 2).Value = "Saksı"
$ws.Cells.Item(39, 3).Value = 982.9
$ws.Cells.Item(40, 1).Value = 38
$ws.Cells.Item(40, 2).Value = "Tohum"
$ws.Cells.Item(40, 3).Value = 804.19
$ws.Cells.Item(41, 1).Value = 39
$ws.Cells.Item(41, 2).Value = "Ajanda"
$ws.Cells.Item(41, 3).Value = 51.28
$ws.Cells.Item(42, 1).Value = 40
$ws.Cells.Item(42, 2).Value = "Boya Seti"
$ws.Cells.Item(42, 3).Value = 41.96
$ws.Cells.Item(43, 1).Value = 41
$ws.Cells.Item(43, 2).Value = "Defter"
$ws.Cells.Item(43, 3).Value = 41.96
$ws.Cells.Item(44, 1).Value = 42
$ws.Cells.Item(44, 2).Value = "Dosya"
$ws.Cells.Item(44, 3).Value = 51.28

